# Banning underscore in user given names (e.g. units and groups):
# rename the "ETS_CO2" group label to "ETS-CO2" throughout the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("emissiondata")

$ws.Range("D2:D5").Value2 = "ETS-CO2"
